# Generate Report for Handback
# Applies the "handback" report-generation edit:
#  - Overview sheet: status text "In Translation" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: fill in "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns (I/J/K) now that handback has happened,
#    with hyperlinks on the new "Latest Target File" cells, and widen a few columns
#    so the new long file names/timestamps are readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("J2").Value = "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.5cd4d04a010a53012ff23f7c577dcbe001780c83.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-21 01:21:30"
$zhcn.Range("J3").Value = "8eb155f8-f613-4999-89a3-692a979094c2.43e81218148976553f3c02e6eead151f5c09aef2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-21 01:21:30"

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md", "", "", "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md", "", "", "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/8eb155f8-f613-4999-89a3-692a979094c2.md", "", "", "8eb155f8-f613-4999-89a3-692a979094c2.md")
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/8eb155f8-f613-4999-89a3-692a979094c2.md", "", "", "8eb155f8-f613-4999-89a3-692a979094c2.md")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("J2").Value = "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.5cd4d04a010a53012ff23f7c577dcbe001780c83.de-de.xlf"
$dede.Range("K2").Value = "2016-10-21 01:21:48"
$dede.Range("J3").Value = "8eb155f8-f613-4999-89a3-692a979094c2.43e81218148976553f3c02e6eead151f5c09aef2.de-de.xlf"
$dede.Range("K3").Value = "2016-10-21 01:21:48"

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md", "", "", "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md")
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md", "", "", "4c17b39e-89ce-4a1f-a9d7-92c3f4b23213.md")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/8eb155f8-f613-4999-89a3-692a979094c2.md", "", "", "8eb155f8-f613-4999-89a3-692a979094c2.md")
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fe3416b2bab57ddd8ee3a3e24dab85bc1c307d7e/e2e/8eb155f8-f613-4999-89a3-692a979094c2.md", "", "", "8eb155f8-f613-4999-89a3-692a979094c2.md")
